{"js": "// Replace each three-digit division problem's text with its new value.\n// The document is a practice worksheet table; every populated cell holds\n// a single run of the form \"ABC\u00f7D=EF, G\". We locate each old expression\n// with Body.search (exact, case-sensitive, not a partial/sub-word match)\n// and overwrite just that run's text, preserving its formatting.\nconst replacements = [\n  [\n    \"196\u00f79=21, 7\",\n    \"321\u00f72=160, 1\"\n  ],\n  [\n    \"947\u00f77=135, 2\",\n    \"536\u00f72=268, 0\"\n  ],\n  [\n    \"226\u00f72=113, 0\",\n    \"517\u00f76=86, 1\"\n  ],\n  [\n    \"370\u00f77=52, 6\",\n    \"275\u00f79=30, 5\"\n  ],\n  [\n    \"916\u00f74=229, 0\",\n    \"870\u00f73=290, 0\"\n  ],\n  [\n    \"914\u00f76=152, 2\",\n    \"637\u00f72=318, 1\"\n  ],\n  [\n    \"229\u00f78=28, 5\",\n    \"963\u00f72=481, 1\"\n  ],\n  [\n    \"681\u00f78=85, 1\",\n    \"119\u00f73=39, 2\"\n  ],\n  [\n    \"382\u00f75=76, 2\",\n    \"372\u00f78=46, 4\"\n  ],\n  [\n    \"231\u00f74=57, 3\",\n    \"699\u00f78=87, 3\"\n  ],\n  [\n    \"576\u00f76=96, 0\",\n    \"654\u00f76=109, 0\"\n  ],\n  [\n    \"334\u00f77=47, 5\",\n    \"242\u00f78=30, 2\"\n  ],\n  [\n    \"559\u00f73=186, 1\",\n    \"236\u00f76=39, 2\"\n  ],\n  [\n    \"271\u00f77=38, 5\",\n    \"397\u00f75=79, 2\"\n  ],\n  [\n    \"249\u00f79=27, 6\",\n    \"495\u00f72=247, 1\"\n  ],\n  [\n    \"846\u00f75=169, 1\",\n    \"963\u00f72=481, 1\"\n  ],\n  [\n    \"570\u00f77=81, 3\",\n    \"879\u00f77=125, 4\"\n  ],\n  [\n    \"972\u00f73=324, 0\",\n    \"555\u00f78=69, 3\"\n  ],\n  [\n    \"341\u00f78=42, 5\",\n    \"638\u00f74=159, 2\"\n  ],\n  [\n    \"172\u00f74=43, 0\",\n    \"450\u00f76=75, 0\"\n  ],\n  [\n    \"606\u00f79=67, 3\",\n    \"109\u00f78=13, 5\"\n  ],\n  [\n    \"634\u00f77=90, 4\",\n    \"874\u00f76=145, 4\"\n  ],\n  [\n    \"993\u00f73=331, 0\",\n    \"114\u00f75=22, 4\"\n  ],\n  [\n    \"936\u00f74=234, 0\",\n    \"393\u00f79=43, 6\"\n  ],\n  [\n    \"863\u00f75=172, 3\",\n    \"128\u00f79=14, 2\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each three-digit division problem's text with its new value.\n# The document is a practice worksheet table; every populated cell holds\n# a single run of the form \"ABC/D=EF, G\" (division sign). We use\n# Range.Find.Execute with MatchCase to locate the exact old expression\n# and replace it in one shot (wdReplaceOne), preserving run formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"196\u00f79=21, 7\"; New = \"321\u00f72=160, 1\" }\n    @{ Old = \"947\u00f77=135, 2\"; New = \"536\u00f72=268, 0\" }\n    @{ Old = \"226\u00f72=113, 0\"; New = \"517\u00f76=86, 1\" }\n    @{ Old = \"370\u00f77=52, 6\"; New = \"275\u00f79=30, 5\" }\n    @{ Old = \"916\u00f74=229, 0\"; New = \"870\u00f73=290, 0\" }\n    @{ Old = \"914\u00f76=152, 2\"; New = \"637\u00f72=318, 1\" }\n    @{ Old = \"229\u00f78=28, 5\"; New = \"963\u00f72=481, 1\" }\n    @{ Old = \"681\u00f78=85, 1\"; New = \"119\u00f73=39, 2\" }\n    @{ Old = \"382\u00f75=76, 2\"; New = \"372\u00f78=46, 4\" }\n    @{ Old = \"231\u00f74=57, 3\"; New = \"699\u00f78=87, 3\" }\n    @{ Old = \"576\u00f76=96, 0\"; New = \"654\u00f76=109, 0\" }\n    @{ Old = \"334\u00f77=47, 5\"; New = \"242\u00f78=30, 2\" }\n    @{ Old = \"559\u00f73=186, 1\"; New = \"236\u00f76=39, 2\" }\n    @{ Old = \"271\u00f77=38, 5\"; New = \"397\u00f75=79, 2\" }\n    @{ Old = \"249\u00f79=27, 6\"; New = \"495\u00f72=247, 1\" }\n    @{ Old = \"846\u00f75=169, 1\"; New = \"963\u00f72=481, 1\" }\n    @{ Old = \"570\u00f77=81, 3\"; New = \"879\u00f77=125, 4\" }\n    @{ Old = \"972\u00f73=324, 0\"; New = \"555\u00f78=69, 3\" }\n    @{ Old = \"341\u00f78=42, 5\"; New = \"638\u00f74=159, 2\" }\n    @{ Old = \"172\u00f74=43, 0\"; New = \"450\u00f76=75, 0\" }\n    @{ Old = \"606\u00f79=67, 3\"; New = \"109\u00f78=13, 5\" }\n    @{ Old = \"634\u00f77=90, 4\"; New = \"874\u00f76=145, 4\" }\n    @{ Old = \"993\u00f73=331, 0\"; New = \"114\u00f75=22, 4\" }\n    @{ Old = \"936\u00f74=234, 0\"; New = \"393\u00f79=43, 6\" }\n    @{ Old = \"863\u00f75=172, 3\"; New = \"128\u00f79=14, 2\" }\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute(\n        $pair.Old,   # FindText\n        $true,       # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceOne)\n    )\n    if (-not $found) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n\n"}
